$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.192991033182664
$ws.Range("C2").Value = 8.56909040750612
$ws.Range("D2").Value = 0.614363312983159

$ws.Range("B3").Value = 0.192730570602027
$ws.Range("C3").Value = 8.56679896154578
$ws.Range("D3").Value = 0.614577468829989

$ws.Range("B4").Value = 0.124594479914413
$ws.Range("C4").Value = 8.80486277888415
$ws.Range("D4").Value = 0.58844822048728

$ws.Range("B5").Value = 0.121972313249735
$ws.Range("C5").Value = 8.83299633761377
$ws.Range("D5").Value = 0.586114801256328

$ws.Range("B6").Value = 0.122147163005381
$ws.Range("C6").Value = 8.82899857825463
$ws.Range("D6").Value = 0.586491437320739

$ws.Range("B7").Value = 0.122421159255852
$ws.Range("C7").Value = 8.82262416091251
$ws.Range("D7").Value = 0.586839200874196

$ws.Range("B8").Value = 0.121776338722119
$ws.Range("C8").Value = 8.83447559959065
$ws.Range("D8").Value = 0.585966306821123

$ws.Range("B9").Value = 0.121658889881121
$ws.Range("C9").Value = 8.83377670241619
$ws.Range("D9").Value = 0.585972914629189

$ws.Range("B10").Value = 0.122959065564311
$ws.Range("C10").Value = 8.82002606903799
$ws.Range("D10").Value = 0.587320875318748
